$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Hspg2"
$ws.Cells.Item(2, 3).Value = "Itga2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 155.2138263333333
$ws.Cells.Item(2, 8).Value = 465.641479
$ws.Cells.Item(2, 9).Value = 0.3492508712612995
$ws.Cells.Item(2, 10).Value = 0.3492508712612995
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.740822
$ws.Cells.Item(2, 14).Value = 5.222466
$ws.Cells.Item(2, 15).Value = 0.4863878955914668
$ws.Cells.Item(2, 16).Value = 0.4863878955914669
$ws.Cells.Item(2, 17).Value = 270.199643585246
$ws.Cells.Item(2, 18).Value = 2431.796792267214
$ws.Cells.Item(2, 19).Value = 0.1698713963062698
$ws.Cells.Item(2, 20).Value = 0.1698713963062698
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Hspg2"
$ws.Cells.Item(3, 3).Value = "Itga2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 155.2138263333333
$ws.Cells.Item(3, 8).Value = 465.641479
$ws.Cells.Item(3, 9).Value = 0.3492508712612995
$ws.Cells.Item(3, 10).Value = 0.3492508712612995
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.077748
$ws.Cells.Item(3, 14).Value = 3.233244
$ws.Cells.Item(3, 15).Value = 0.3011241710513264
$ws.Cells.Item(3, 16).Value = 0.3011241710513265
$ws.Cells.Item(3, 17).Value = 167.2813909030973
$ws.Cells.Item(3, 18).Value = 1505.532518127876
$ws.Cells.Item(3, 19).Value = 0.1051678790975123
$ws.Cells.Item(3, 20).Value = 0.1051678790975124
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Hspg2"
$ws.Cells.Item(4, 3).Value = "Itga2"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 155.2138263333333
$ws.Cells.Item(4, 8).Value = 465.641479
$ws.Cells.Item(4, 9).Value = 0.3492508712612995
$ws.Cells.Item(4, 10).Value = 0.3492508712612995
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.03488166666666666
$ws.Cells.Item(4, 14).Value = 0.104645
$ws.Cells.Item(4, 15).Value = 0.009745982326006345
$ws.Cells.Item(4, 16).Value = 0.009745982326006345
$ws.Cells.Item(4, 17).Value = 5.414116952217222
$ws.Cells.Item(4, 18).Value = 48.727052569955
$ws.Cells.Item(4, 19).Value = 0.003403792818654943
$ws.Cells.Item(4, 20).Value = 0.003403792818654943
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Hspg2"
$ws.Cells.Item(5, 3).Value = "Itga2"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 155.2138263333333
$ws.Cells.Item(5, 8).Value = 465.641479
$ws.Cells.Item(5, 9).Value = 0.3492508712612995
$ws.Cells.Item(5, 10).Value = 0.3492508712612995
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.7256300000000001
$ws.Cells.Item(5, 14).Value = 2.17689
$ws.Cells.Item(5, 15).Value = 0.2027419510312003
$ws.Cells.Item(5, 16).Value = 0.2027419510312003
$ws.Cells.Item(5, 17).Value = 112.6278088022567
$ws.Cells.Item(5, 18).Value = 1013.65027922031
$ws.Cells.Item(5, 19).Value = 0.07080780303886242
$ws.Cells.Item(5, 20).Value = 0.07080780303886242
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Hspg2"
$ws.Cells.Item(6, 3).Value = "Itga2"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 246.1811623333333
$ws.Cells.Item(6, 8).Value = 738.543487
$ws.Cells.Item(6, 9).Value = 0.5539389593320749
$ws.Cells.Item(6, 10).Value = 0.5539389593320749
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.740822
$ws.Cells.Item(6, 14).Value = 5.222466
$ws.Cells.Item(6, 15).Value = 0.4863878955914668
$ws.Cells.Item(6, 16).Value = 0.4863878955914669
$ws.Cells.Item(6, 17).Value = 428.557583375438
$ws.Cells.Item(6, 18).Value = 3857.018250378942
$ws.Cells.Item(6, 19).Value = 0.2694292047156551
$ws.Cells.Item(6, 20).Value = 0.2694292047156551
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Hspg2"
$ws.Cells.Item(7, 3).Value = "Itga2"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 246.1811623333333
$ws.Cells.Item(7, 8).Value = 738.543487
$ws.Cells.Item(7, 9).Value = 0.5539389593320749
$ws.Cells.Item(7, 10).Value = 0.5539389593320749
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.077748
$ws.Cells.Item(7, 14).Value = 3.233244
$ws.Cells.Item(7, 15).Value = 0.3011241710513264
$ws.Cells.Item(7, 16).Value = 0.3011241710513265
$ws.Cells.Item(7, 17).Value = 265.3212553424253
$ws.Cells.Item(7, 18).Value = 2387.891298081828
$ws.Cells.Item(7, 19).Value = 0.1668044099419055
$ws.Cells.Item(7, 20).Value = 0.1668044099419055
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Hspg2"
$ws.Cells.Item(8, 3).Value = "Itga2"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 246.1811623333333
$ws.Cells.Item(8, 8).Value = 738.543487
$ws.Cells.Item(8, 9).Value = 0.5539389593320749
$ws.Cells.Item(8, 10).Value = 0.5539389593320749
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.03488166666666666
$ws.Cells.Item(8, 14).Value = 0.104645
$ws.Cells.Item(8, 15).Value = 0.009745982326006345
$ws.Cells.Item(8, 16).Value = 0.009745982326006345
$ws.Cells.Item(8, 17).Value = 8.587209244123889
$ws.Cells.Item(8, 18).Value = 77.28488319711501
$ws.Cells.Item(8, 19).Value = 0.00539867930733675
$ws.Cells.Item(8, 20).Value = 0.00539867930733675
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Hspg2"
$ws.Cells.Item(9, 3).Value = "Itga2"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 246.1811623333333
$ws.Cells.Item(9, 8).Value = 738.543487
$ws.Cells.Item(9, 9).Value = 0.5539389593320749
$ws.Cells.Item(9, 10).Value = 0.5539389593320749
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.7256300000000001
$ws.Cells.Item(9, 14).Value = 2.17689
$ws.Cells.Item(9, 15).Value = 0.2027419510312003
$ws.Cells.Item(9, 16).Value = 0.2027419510312003
$ws.Cells.Item(9, 17).Value = 178.6364368239367
$ws.Cells.Item(9, 18).Value = 1607.72793141543
$ws.Cells.Item(9, 19).Value = 0.1123066653671776
$ws.Cells.Item(9, 20).Value = 0.1123066653671776
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Hspg2"
$ws.Cells.Item(10, 3).Value = "Itga2"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.2401933333333333
$ws.Cells.Item(10, 8).Value = 0.72058
$ws.Cells.Item(10, 9).Value = 0.0005404655817044752
$ws.Cells.Item(10, 10).Value = 0.0005404655817044752
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.740822
$ws.Cells.Item(10, 14).Value = 5.222466
$ws.Cells.Item(10, 15).Value = 0.4863878955914668
$ws.Cells.Item(10, 16).Value = 0.4863878955914669
$ws.Cells.Item(10, 17).Value = 0.41813383892
$ws.Cells.Item(10, 18).Value = 3.76320455028
$ws.Cells.Item(10, 19).Value = 0.0002628759169248577
$ws.Cells.Item(10, 20).Value = 0.0002628759169248577
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Hspg2"
$ws.Cells.Item(11, 3).Value = "Itga2"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.2401933333333333
$ws.Cells.Item(11, 8).Value = 0.72058
$ws.Cells.Item(11, 9).Value = 0.0005404655817044752
$ws.Cells.Item(11, 10).Value = 0.0005404655817044752
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.077748
$ws.Cells.Item(11, 14).Value = 3.233244
$ws.Cells.Item(11, 15).Value = 0.3011241710513264
$ws.Cells.Item(11, 16).Value = 0.3011241710513265
$ws.Cells.Item(11, 17).Value = 0.2588678846133333
$ws.Cells.Item(11, 18).Value = 2.32981096152
$ws.Cells.Item(11, 19).Value = 0.000162747250272533
$ws.Cells.Item(11, 20).Value = 0.0001627472502725331
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Hspg2"
$ws.Cells.Item(12, 3).Value = "Itga2"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.2401933333333333
$ws.Cells.Item(12, 8).Value = 0.72058
$ws.Cells.Item(12, 9).Value = 0.0005404655817044752
$ws.Cells.Item(12, 10).Value = 0.0005404655817044752
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.03488166666666666
$ws.Cells.Item(12, 14).Value = 0.104645
$ws.Cells.Item(12, 15).Value = 0.009745982326006345
$ws.Cells.Item(12, 16).Value = 0.009745982326006345
$ws.Cells.Item(12, 17).Value = 0.008378343788888888
$ws.Cells.Item(12, 18).Value = 0.0754050941
$ws.Cells.Item(12, 19).Value = 0.000005267368007106554
$ws.Cells.Item(12, 20).Value = 0.000005267368007106554
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Hspg2"
$ws.Cells.Item(13, 3).Value = "Itga2"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.2401933333333333
$ws.Cells.Item(13, 8).Value = 0.72058
$ws.Cells.Item(13, 9).Value = 0.0005404655817044752
$ws.Cells.Item(13, 10).Value = 0.0005404655817044752
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.7256300000000001
$ws.Cells.Item(13, 14).Value = 2.17689
$ws.Cells.Item(13, 15).Value = 0.2027419510312003
$ws.Cells.Item(13, 16).Value = 0.2027419510312003
$ws.Cells.Item(13, 17).Value = 0.1742914884666667
$ws.Cells.Item(13, 18).Value = 1.5686233962
$ws.Cells.Item(13, 19).Value = 0.0001095750464999779
$ws.Cells.Item(13, 20).Value = 0.0001095750464999779
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Hspg2"
$ws.Cells.Item(14, 3).Value = "Itga2"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 42.784114
$ws.Cells.Item(14, 8).Value = 128.352342
$ws.Cells.Item(14, 9).Value = 0.09626970382492123
$ws.Cells.Item(14, 10).Value = 0.09626970382492124
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.740822
$ws.Cells.Item(14, 14).Value = 5.222466
$ws.Cells.Item(14, 15).Value = 0.4863878955914668
$ws.Cells.Item(14, 16).Value = 0.4863878955914669
$ws.Cells.Item(14, 17).Value = 74.47952690170798
$ws.Cells.Item(14, 18).Value = 670.3157421153719
$ws.Cells.Item(14, 19).Value = 0.04682441865261722
$ws.Cells.Item(14, 20).Value = 0.04682441865261723
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Hspg2"
$ws.Cells.Item(15, 3).Value = "Itga2"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 42.784114
$ws.Cells.Item(15, 8).Value = 128.352342
$ws.Cells.Item(15, 9).Value = 0.09626970382492123
$ws.Cells.Item(15, 10).Value = 0.09626970382492124
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 1.077748
$ws.Cells.Item(15, 14).Value = 3.233244
$ws.Cells.Item(15, 15).Value = 0.3011241710513264
$ws.Cells.Item(15, 16).Value = 0.3011241710513265
$ws.Cells.Item(15, 17).Value = 46.11049329527199
$ws.Cells.Item(15, 18).Value = 414.994439657448
$ws.Cells.Item(15, 19).Value = 0.02898913476163612
$ws.Cells.Item(15, 20).Value = 0.02898913476163613
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Hspg2"
$ws.Cells.Item(16, 3).Value = "Itga2"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 42.784114
$ws.Cells.Item(16, 8).Value = 128.352342
$ws.Cells.Item(16, 9).Value = 0.09626970382492123
$ws.Cells.Item(16, 10).Value = 0.09626970382492124
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.03488166666666666
$ws.Cells.Item(16, 14).Value = 0.104645
$ws.Cells.Item(16, 15).Value = 0.009745982326006345
$ws.Cells.Item(16, 16).Value = 0.009745982326006345
$ws.Cells.Item(16, 17).Value = 1.492381203176666
$ws.Cells.Item(16, 18).Value = 13.43143082859
$ws.Cells.Item(16, 19).Value = 0.0009382428320075478
$ws.Cells.Item(16, 20).Value = 0.0009382428320075479
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Hspg2"
$ws.Cells.Item(17, 3).Value = "Itga2"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 42.784114
$ws.Cells.Item(17, 8).Value = 128.352342
$ws.Cells.Item(17, 9).Value = 0.09626970382492123
$ws.Cells.Item(17, 10).Value = 0.09626970382492124
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.7256300000000001
$ws.Cells.Item(17, 14).Value = 2.17689
$ws.Cells.Item(17, 15).Value = 0.2027419510312003
$ws.Cells.Item(17, 16).Value = 0.2027419510312003
$ws.Cells.Item(17, 17).Value = 31.04543664182
$ws.Cells.Item(17, 18).Value = 279.40892977638
$ws.Cells.Item(17, 19).Value = 0.01951790757866034
$ws.Cells.Item(17, 20).Value = 0.01951790757866034
